# Append new incident rows (127-138) to the management log on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("WC47 NACP",  "Fallo en elevador",      "2024-06-04", "09:18:31", "Mañana", "09:18:32", "0:00:01", "N/A"),
    @("WC47 NACP",  "Ascensor no sube",       "2024-06-04", "09:19:05", "Mañana", "09:19:06", "0:00:01", "N/A"),
    @("WC47 NACP",  "Fallo tornillo",         "2024-06-04", "09:19:18", "Mañana", "09:19:19", "0:00:01", "0.57 minutos"),
    @("WC47 NACP",  "Fallo tornillo",         "2024-06-04", "09:19:57", "Mañana", "09:19:57", "0:00:00", "0.39 minutos"),
    @("WC47 NACP",  "Fallo en elevador",      "2024-06-04", "09:20:17", "Mañana", "09:20:18", "0:00:01", "0.48 minutos"),
    @("WC47 NACP",  "No coge placa",          "2024-06-04", "09:21:12", "Mañana", "09:21:13", "0:00:01", "N/A"),
    @("WC47 NACP",  "Etiquetadora",           "2024-06-04", "09:26:54", "Mañana", "09:26:55", "0:00:01", "N/A"),
    @("WC47 NACP",  "Etiquetadora",           "2024-06-04", "09:26:56", "Mañana", "09:26:57", "0:00:01", "5.70 minutos"),
    @("WC47 NACP",  "No atornilla tapa",      "2024-06-04", "09:36:19", "Mañana", "09:36:21", "0:00:02", "2.87 minutos"),
    @("WC47 NACP",  "Fallo en paletizador",   "2024-06-04", "09:36:23", "Mañana", "09:36:24", "0:00:01", "5.04 minutos"),
    @("WC47 NACP",  "Fallo en paletizador",   "2024-06-04", "09:36:25", "Mañana", "09:36:27", "0:00:02", "3.80 minutos"),
    @("WV50 FILTER","Fallo cámara ferrite",   "2024-06-04", "09:52:34", "Mañana", "09:52:35", "0:00:01", "N/A")
)

$startRow = 127
$endRow = $startRow + $rows.Count - 1

# Pre-format the target block as Text so date/time-looking strings (e.g.
# "2024-06-04", "09:18:31") are stored verbatim instead of being parsed
# into Excel date/time serials, matching the source data's plain-text cells.
$blockRange = $ws.Range("A$startRow`:H$endRow")
$blockRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Drop back to the default (General/Normal) style now that the text values
# are committed, so the new cells carry no explicit style index.
$blockRange.Style = "Normal"
